$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("S2").Value = 2.82
$ws.Range("T2").Value = 1.64
$ws.Range("U2").Value = 2.44
$ws.Range("Y2").Value = 16.5
$ws.Range("AH2").Value = 16
$ws.Range("AJ2").Value = 29
$ws.Range("AO2").Value = 27
$ws.Range("S3").Value = 3.05
$ws.Range("G4").Value = 2.12
$ws.Range("J4").Value = 3.25
$ws.Range("W4").Value = 1.89
$ws.Range("AD4").Value = 1000
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 5.1
$ws.Range("O5").Value = 1.28
$ws.Range("R5").Value = 1.38
$ws.Range("S5").Value = 3.15
$ws.Range("I6").Value = 1.44
$ws.Range("V6").Value = 3.2
$ws.Range("Z6").Value = 9
$ws.Range("Q7").Value = 1.36
$ws.Range("AN7").Value = 3.35
$ws.Range("U8").Value = 2.02
$ws.Range("K12").Value = 5.3
$ws.Range("P12").Value = 3.2
$ws.Range("T12").Value = 1.42
$ws.Range("U12").Value = 2.84
$ws.Range("V12").Value = 1.3
$ws.Range("W12").Value = 1.95
$ws.Range("L13").Value = 1.23
$ws.Range("R13").Value = 1.74
$ws.Range("G14").Value = 2.2
$ws.Range("I14").Value = 5
$ws.Range("J14").Value = 3.5
$ws.Range("R14").Value = 1.48
$ws.Range("S14").Value = 1.59
$ws.Range("V14").Value = 1.25
$ws.Range("W14").Value = 1.83
$ws.Range("V16").Value = 1.77
$ws.Range("J17").Value = 5.2
$ws.Range("M17").Value = 1.02
$ws.Range("N17").Value = 6.6
$ws.Range("O17").Value = 1.14
$ws.Range("Q17").Value = 1.39
$ws.Range("R17").Value = 1.76
$ws.Range("S17").Value = 2.06
$ws.Range("T17").Value = 1.64
$ws.Range("U17").Value = 2.26
$ws.Range("W17").Value = 2.96
$ws.Range("X17").Value = 40
$ws.Range("Y17").Value = 42
$ws.Range("Z17").Value = 80
$ws.Range("AA17").Value = 210
$ws.Range("AB17").Value = 16
$ws.Range("AC17").Value = 15.5
$ws.Range("AD17").Value = 32
$ws.Range("AE17").Value = 95
$ws.Range("AF17").Value = 14
$ws.Range("AG17").Value = 12.5
$ws.Range("AH17").Value = 24
$ws.Range("AI17").Value = 75
$ws.Range("AJ17").Value = 16.5
$ws.Range("AK17").Value = 16.5
$ws.Range("AL17").Value = 30
$ws.Range("AM17").Value = 85
$ws.Range("AN17").Value = 5.2
$ws.Range("AO17").Value = 75
$ws.Range("F19").Value = 1.59
$ws.Range("G20").Value = 2.94
$ws.Range("I20").Value = 4.7
$ws.Range("J20").Value = 2.6
$ws.Range("V20").Value = 1.27
$ws.Range("W20").Value = 1.51
$ws.Range("J21").Value = 3.5
$ws.Range("N22").Value = 1.64
$ws.Range("P22").Value = 1.64
$ws.Range("L23").Value = 1.58
$ws.Range("O23").Value = 1.53
$ws.Range("Q23").Value = 2.68
$ws.Range("U23").Value = 1.84
$ws.Range("F24").Value = 2.3
$ws.Range("I24").Value = 4.4
$ws.Range("J24").Value = 2.24
$ws.Range("K24").Value = 3.8
$ws.Range("O24").Value = 1.67
$ws.Range("V24").Value = 1.29
$ws.Range("F26").Value = 2.16
$ws.Range("P27").Value = 1.78
$ws.Range("S27").Value = 3.45
$ws.Range("AF27").Value = 8
$ws.Range("AJ27").Value = 12.5
$ws.Range("G28").Value = 3.2
$ws.Range("I28").Value = 3.35
$ws.Range("L28").Value = 1.01
$ws.Range("M28").Value = 1.01
$ws.Range("N28").Value = 1.6
$ws.Range("O28").Value = 1.01
$ws.Range("Q28").Value = 1.47
$ws.Range("R28").Value = 1.08
$ws.Range("S28").Value = 1.01
$ws.Range("T28").Value = 1.01
$ws.Range("U28").Value = 1.01
$ws.Range("V28").Value = 1.42
$ws.Range("W28").Value = 1.46
$ws.Range("X28").Value = 1000
$ws.Range("Y28").Value = 1000
$ws.Range("Z28").Value = 1000
$ws.Range("AA28").Value = 1000
$ws.Range("AB28").Value = 1000
$ws.Range("AC28").Value = 1000
$ws.Range("AD28").Value = 1000
$ws.Range("AE28").Value = 1000
$ws.Range("AF28").Value = 1000
$ws.Range("AG28").Value = 1000
$ws.Range("AH28").Value = 1000
$ws.Range("AI28").Value = 1000
$ws.Range("AJ28").Value = 1000
$ws.Range("AK28").Value = 1000
$ws.Range("AL28").Value = 1000
$ws.Range("AM28").Value = 1000
$ws.Range("AN28").Value = 1000
$ws.Range("AO28").Value = 1000
